$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2649.1667
$ws.Range("J88").Value = 2649.1667
$ws.Range("L88").Value = 2649.1667
$ws.Range("N88").Value = -3461.1667
$ws.Range("H91").Value = 2649.1667
$ws.Range("J91").Value = 2649.1667
$ws.Range("L91").Value = 2649.1667
$ws.Range("N91").Value = -5457.1667
$ws.Range("H94").Value = 2758.4285
$ws.Range("I94").Value = 2758.4285
$ws.Range("K94").Value = 2758.4285
$ws.Range("M94").Value = -2307.4285
$ws.Range("H106").Value = 2789.2727
$ws.Range("I106").Value = 3686.1428
$ws.Range("K106").Value = 3686.1428
$ws.Range("M106").Value = -3055.1428
$ws.Range("H116").Value = 9665.462
$ws.Range("I116").Value = 14377.75
$ws.Range("J116").Value = 2125.8
$ws.Range("K116").Value = 14377.75
$ws.Range("L116").Value = 2125.8
$ws.Range("M116").Value = -10935.75
$ws.Range("N116").Value = -9009.799999999999
$ws.Range("H125").Value = 555.55554
$ws.Range("I125").Value = 555.55554
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4999.99986
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2539.99986
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1377
$ws.Range("I132").Value = 1344.1538
$ws.Range("K132").Value = 4032.4614
$ws.Range("M132").Value = -1502.4614
$ws.Range("H135").Value = 619.4
$ws.Range("I135").Value = 242
$ws.Range("J135").Value = 1185.5
$ws.Range("K135").Value = 2178
$ws.Range("L135").Value = 10669.5
$ws.Range("M135").Value = 357
$ws.Range("N135").Value = -15739.5
$ws.Range("H137").Value = 1356.2778
$ws.Range("I137").Value = 1090.7
$ws.Range("K137").Value = 3272.1
$ws.Range("M137").Value = -722.1000000000004
$ws.Range("H138").Value = 3175.5112
$ws.Range("I138").Value = 4657.6
$ws.Range("K138").Value = 13972.8
$ws.Range("M138").Value = -8832.800000000001
$ws.Range("H141").Value = 2195.7144
$ws.Range("I141").Value = 853.2
$ws.Range("J141").Value = 5552
$ws.Range("K141").Value = 2559.6
$ws.Range("L141").Value = 16656
$ws.Range("M141").Value = 2620.4
$ws.Range("N141").Value = -27016

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3927.897
$ws.Range("I32").Value = 3470.8154
$ws.Range("K32").Value = 3470.8154
$ws.Range("M32").Value = -3183.8154
$ws.Range("H54").Value = 17320
$ws.Range("J54").Value = 17320
$ws.Range("L54").Value = 17320
$ws.Range("N54").Value = -18858
$ws.Range("H61").Value = 5552.44
$ws.Range("I61").Value = 5495.6
$ws.Range("K61").Value = 5495.6
$ws.Range("M61").Value = -5283.6
$ws.Range("H74").Value = 1106.2142
$ws.Range("I74").Value = 457.25
$ws.Range("K74").Value = 457.25
$ws.Range("M74").Value = 416.75
$ws.Range("H77").Value = 1106.2142
$ws.Range("I77").Value = 457.25
$ws.Range("K77").Value = 2286.25
$ws.Range("M77").Value = 2081.75
$ws.Range("H136").Value = 5552.44
$ws.Range("I136").Value = 5495.6
$ws.Range("K136").Value = 16486.8
$ws.Range("M136").Value = -13936.8

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2379
$ws.Range("I20").Value = 1859.125
$ws.Range("J20").Value = 3418.75
$ws.Range("K20").Value = 1859.125
$ws.Range("L20").Value = 3418.75
$ws.Range("M20").Value = -1612.125
$ws.Range("N20").Value = -3912.75
$ws.Range("H94").Value = 2007.2667
$ws.Range("I94").Value = 512.7
$ws.Range("J94").Value = 4996.4
$ws.Range("K94").Value = 512.7
$ws.Range("L94").Value = 4996.4
$ws.Range("M94").Value = -61.70000000000005
$ws.Range("N94").Value = -5898.4
$ws.Range("H99").Value = 1524.7
$ws.Range("J99").Value = 1873.75
$ws.Range("L99").Value = 1873.75
$ws.Range("N99").Value = -4869.75
$ws.Range("H108").Value = 65000
$ws.Range("J108").Value = 65000
$ws.Range("L108").Value = 65000
$ws.Range("N108").Value = -72680
$ws.Range("H134").Value = 4546.7812
$ws.Range("I134").Value = 4776.6787
$ws.Range("K134").Value = 14330.0361
$ws.Range("M134").Value = -11795.0361

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4250
$ws.Range("I6").Value = 4250
$ws.Range("K6").Value = 4250
$ws.Range("M6").Value = -4137
$ws.Range("H16").Value = 598.8946999999999
$ws.Range("I16").Value = 533.46155
$ws.Range("K16").Value = 533.46155
$ws.Range("M16").Value = -246.46155
$ws.Range("H31").Value = 2558.3635
$ws.Range("I31").Value = 2200
$ws.Range("J31").Value = 2988.4
$ws.Range("K31").Value = 2200
$ws.Range("L31").Value = 2988.4
$ws.Range("M31").Value = -1905
$ws.Range("N31").Value = -3578.4
$ws.Range("H34").Value = 2558.3635
$ws.Range("I34").Value = 2200
$ws.Range("J34").Value = 2988.4
$ws.Range("K34").Value = 2200
$ws.Range("L34").Value = 2988.4
$ws.Range("M34").Value = -1998
$ws.Range("N34").Value = -3392.4
$ws.Range("H58").Value = 2072446.4
$ws.Range("I58").Value = 2290177.5
$ws.Range("K58").Value = 2290177.5
$ws.Range("M58").Value = -2289974.5
$ws.Range("H113").Value = 598.8946999999999
$ws.Range("I113").Value = 533.46155
$ws.Range("K113").Value = 533.46155
$ws.Range("M113").Value = 1636.53845
$ws.Range("H132").Value = 1661.15
$ws.Range("I132").Value = 1130.8823
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 3392.6469
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -862.6468999999997
$ws.Range("N132").Value = -19058
$ws.Range("H134").Value = 1841.2
$ws.Range("I134").Value = 1613.4242
$ws.Range("K134").Value = 4840.2726
$ws.Range("M134").Value = -2305.2726
$ws.Range("H136").Value = 2072446.4
$ws.Range("I136").Value = 2290177.5
$ws.Range("K136").Value = 6870532.5
$ws.Range("M136").Value = -6867982.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1012
$ws.Range("I63").Value = 1012
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3036
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2287
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1012
$ws.Range("I66").Value = 1012
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9108
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5364
$ws.Range("N66").ClearContents()
$ws.Range("H98").Value = 1064.6
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H128").Value = 399999
$ws.Range("I128").Value = 399999
$ws.Range("K128").Value = 1199997
$ws.Range("M128").Value = -1195017
$ws.Range("H137").Value = 3398.8076
$ws.Range("I137").Value = 1687.75
$ws.Range("J137").Value = 6136.5
$ws.Range("K137").Value = 5063.25
$ws.Range("L137").Value = 18409.5
$ws.Range("M137").Value = 36.75
$ws.Range("N137").Value = -28609.5
$ws.Range("H140").Value = 1924.4667
$ws.Range("I140").Value = 1055.75
$ws.Range("J140").Value = 5399.3335
$ws.Range("K140").Value = 3167.25
$ws.Range("L140").Value = 16198.0005
$ws.Range("M140").Value = 2012.75
$ws.Range("N140").Value = -26558.0005

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2431
$ws.Range("I102").Value = 2075.4
$ws.Range("K102").Value = 2075.4
$ws.Range("M102").Value = -453.4000000000001
$ws.Range("H113").Value = 1015.0909
$ws.Range("I113").Value = 580.7143
$ws.Range("K113").Value = 580.7143
$ws.Range("M113").Value = 1589.2857
$ws.Range("H126").Value = 2359544.5
$ws.Range("I126").Value = 3706881.2
$ws.Range("J126").Value = 113983.22
$ws.Range("K126").Value = 11120643.6
$ws.Range("L126").Value = 341949.66
$ws.Range("M126").Value = -11118173.6
$ws.Range("N126").Value = -346889.66

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2490.7334
$ws.Range("I7").Value = 2646.2
$ws.Range("J7").Value = 2179.8
$ws.Range("K7").Value = 2646.2
$ws.Range("L7").Value = 2179.8
$ws.Range("M7").Value = -2534.2
$ws.Range("N7").Value = -2403.8
$ws.Range("H126").Value = 2490.7334
$ws.Range("I126").Value = 2646.2
$ws.Range("J126").Value = 2179.8
$ws.Range("K126").Value = 7938.599999999999
$ws.Range("L126").Value = 6539.400000000001
$ws.Range("M126").Value = -5468.599999999999
$ws.Range("N126").Value = -11479.4
$ws.Range("H136").Value = 1955.9474
$ws.Range("I136").Value = 1597
$ws.Range("K136").Value = 4791
$ws.Range("M136").Value = -2241

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H96").Value = 1118.2858
$ws.Range("I96").Value = 499.75
$ws.Range("J96").Value = 1365.7
$ws.Range("K96").Value = 499.75
$ws.Range("L96").Value = 1365.7
$ws.Range("M96").Value = 873.25
$ws.Range("N96").Value = -4111.7
$ws.Range("H126").Value = 1228.5625
$ws.Range("I126").Value = 970.913
$ws.Range("J126").Value = 1887
$ws.Range("K126").Value = 2912.739
$ws.Range("L126").Value = 5661
$ws.Range("M126").Value = -442.739
$ws.Range("N126").Value = -10601
$ws.Range("H132").Value = 1891.2963
$ws.Range("I132").Value = 1364.9546
$ws.Range("J132").Value = 4207.2
$ws.Range("K132").Value = 4094.8638
$ws.Range("L132").Value = 12621.6
$ws.Range("M132").Value = -1564.8638
$ws.Range("N132").Value = -17681.6
